$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.017.03"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.585.84"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "2.598.26"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "3.050.04"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "58.969.59"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "2.593.75"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "0.0₃0724"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.77%  "
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.601"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0951"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0515"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "1.968.24"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.56%  "
